$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 76
$ws_ALC.Range("H76").Value = 3350
$ws_ALC.Range("I76").Value = 3350
$ws_ALC.Range("K76").Value = 3350
$ws_ALC.Range("M76").Value = -3035

# ALC row 79
$ws_ALC.Range("H79").Value = 3350
$ws_ALC.Range("I79").Value = 3350
$ws_ALC.Range("K79").Value = 3350
$ws_ALC.Range("M79").Value = -2258

# ALC row 100
$ws_ALC.Range("H100").Value = 16667946
$ws_ALC.Range("I100").Value = 20001346
$ws_ALC.Range("J100").Value = 950
$ws_ALC.Range("K100").Value = 20001346
$ws_ALC.Range("L100").Value = 950
$ws_ALC.Range("M100").Value = -20000805
$ws_ALC.Range("N100").Value = -2032

# ALC row 113
$ws_ALC.Range("H113").Value = 10227.818
$ws_ALC.Range("I113").Value = 4300
$ws_ALC.Range("K113").Value = 4300
$ws_ALC.Range("M113").Value = -1046

# ALC row 124
$ws_ALC.Range("H124").Value = 42446.668
$ws_ALC.Range("J124").Value = 42446.668
$ws_ALC.Range("L124").Value = 42446.668
$ws_ALC.Range("N124").Value = -52266.668

# ALC row 132
$ws_ALC.Range("H132").Value = 111419.54
$ws_ALC.Range("I132").Value = 135484.44
$ws_ALC.Range("J132").Value = 8284.286
$ws_ALC.Range("K132").Value = 406453.32
$ws_ALC.Range("L132").Value = 24852.858
$ws_ALC.Range("M132").Value = -403923.32
$ws_ALC.Range("N132").Value = -29912.858

# ALC row 141
$ws_ALC.Range("H141").Value = 8978.666999999999
$ws_ALC.Range("I141").Value = 9590.77
$ws_ALC.Range("K141").Value = 28772.31
$ws_ALC.Range("M141").Value = -23592.31

# ARM row 2
$ws_ARM.Range("H2").Value = 958.4167
$ws_ARM.Range("I2").Value = 1082.875
$ws_ARM.Range("J2").Value = 709.5
$ws_ARM.Range("K2").Value = 1082.875
$ws_ARM.Range("L2").Value = 709.5
$ws_ARM.Range("M2").Value = -969.875
$ws_ARM.Range("N2").Value = -935.5

# ARM row 3
$ws_ARM.Range("H3").Value = 8764.4
$ws_ARM.Range("I3").Value = 2502.5
$ws_ARM.Range("J3").Value = 12939
$ws_ARM.Range("K3").Value = 2502.5
$ws_ARM.Range("L3").Value = 12939
$ws_ARM.Range("M3").Value = -2387.5
$ws_ARM.Range("N3").Value = -13169

# ARM row 32
$ws_ARM.Range("H32").Value = 7784.5845
$ws_ARM.Range("I32").Value = 4784.4185
$ws_ARM.Range("J32").Value = 11578.912
$ws_ARM.Range("K32").Value = 4784.4185
$ws_ARM.Range("L32").Value = 11578.912
$ws_ARM.Range("M32").Value = -4497.4185
$ws_ARM.Range("N32").Value = -12152.912

# ARM row 61
$ws_ARM.Range("H61").Value = 1931.2693
$ws_ARM.Range("I61").Value = 1069.6875
$ws_ARM.Range("K61").Value = 1069.6875
$ws_ARM.Range("M61").Value = -857.6875

# ARM row 74
$ws_ARM.Range("H74").Value = 3711.8948
$ws_ARM.Range("I74").Value = 3787.4666
$ws_ARM.Range("J74").Value = 3428.5
$ws_ARM.Range("K74").Value = 3787.4666
$ws_ARM.Range("L74").Value = 3428.5
$ws_ARM.Range("M74").Value = -2913.4666
$ws_ARM.Range("N74").Value = -5176.5

# ARM row 77
$ws_ARM.Range("H77").Value = 3711.8948
$ws_ARM.Range("I77").Value = 3787.4666
$ws_ARM.Range("J77").Value = 3428.5
$ws_ARM.Range("K77").Value = 18937.333
$ws_ARM.Range("L77").Value = 17142.5
$ws_ARM.Range("M77").Value = -14569.333
$ws_ARM.Range("N77").Value = -25878.5

# ARM row 110
$ws_ARM.Range("H110").Value = 1053.6364
$ws_ARM.Range("I110").Value = 981.6667
$ws_ARM.Range("J110").Value = 1140
$ws_ARM.Range("K110").Value = 981.6667
$ws_ARM.Range("L110").Value = 1140
$ws_ARM.Range("M110").Value = 1063.3333
$ws_ARM.Range("N110").Value = -5230

# ARM row 116
$ws_ARM.Range("H116").Value = 958.4167
$ws_ARM.Range("I116").Value = 1082.875
$ws_ARM.Range("J116").Value = 709.5
$ws_ARM.Range("K116").Value = 1082.875
$ws_ARM.Range("L116").Value = 709.5
$ws_ARM.Range("M116").Value = 1211.125
$ws_ARM.Range("N116").Value = -5297.5

# ARM row 136
$ws_ARM.Range("H136").Value = 1931.2693
$ws_ARM.Range("I136").Value = 1069.6875
$ws_ARM.Range("K136").Value = 3209.0625
$ws_ARM.Range("M136").Value = -659.0625

# BSM row 3
$ws_BSM.Range("H3").Value = 958.4167
$ws_BSM.Range("I3").Value = 1082.875
$ws_BSM.Range("J3").Value = 709.5
$ws_BSM.Range("K3").Value = 1082.875
$ws_BSM.Range("L3").Value = 709.5
$ws_BSM.Range("M3").Value = -968.875
$ws_BSM.Range("N3").Value = -937.5

# BSM row 86
$ws_BSM.Range("H86").Value = 2569
$ws_BSM.Range("I86").Value = 2850
$ws_BSM.Range("J86").Value = 2007
$ws_BSM.Range("K86").Value = 2850
$ws_BSM.Range("L86").Value = 2007
$ws_BSM.Range("M86").Value = -1727
$ws_BSM.Range("N86").Value = -4253

# BSM row 89
$ws_BSM.Range("H89").Value = 2569
$ws_BSM.Range("I89").Value = 2850
$ws_BSM.Range("J89").Value = 2007
$ws_BSM.Range("K89").Value = 14250
$ws_BSM.Range("L89").Value = 10035
$ws_BSM.Range("M89").Value = -8634
$ws_BSM.Range("N89").Value = -21267

# BSM row 99
$ws_BSM.Range("H99").Value = 2764.8386
$ws_BSM.Range("I99").Value = 1337.3684
$ws_BSM.Range("J99").Value = 5025
$ws_BSM.Range("K99").Value = 1337.3684
$ws_BSM.Range("L99").Value = 5025
$ws_BSM.Range("M99").Value = 160.6315999999999
$ws_BSM.Range("N99").Value = -8021

# BSM row 105
$ws_BSM.Range("H105").Value = 2580.6584
$ws_BSM.Range("I105").Value = 2564.5386
$ws_BSM.Range("K105").Value = 2564.5386
$ws_BSM.Range("M105").Value = -817.5385999999999

# BSM row 134
$ws_BSM.Range("H134").Value = 3676.342
$ws_BSM.Range("I134").Value = 2078.7307
$ws_BSM.Range("K134").Value = 6236.1921
$ws_BSM.Range("M134").Value = -3701.1921

# CRP row 31
$ws_CRP.Range("H31").Value = 2967.8936
$ws_CRP.Range("I31").Value = 1318.8096
$ws_CRP.Range("J31").Value = 4299.846
$ws_CRP.Range("K31").Value = 1318.8096
$ws_CRP.Range("L31").Value = 4299.846
$ws_CRP.Range("M31").Value = -1023.8096
$ws_CRP.Range("N31").Value = -4889.846

# CRP row 34
$ws_CRP.Range("H34").Value = 2967.8936
$ws_CRP.Range("I34").Value = 1318.8096
$ws_CRP.Range("J34").Value = 4299.846
$ws_CRP.Range("K34").Value = 1318.8096
$ws_CRP.Range("L34").Value = 4299.846
$ws_CRP.Range("M34").Value = -1116.8096
$ws_CRP.Range("N34").Value = -4703.846

# CRP row 58
$ws_CRP.Range("H58").Value = 1856.0555
$ws_CRP.Range("I58").Value = 1606.0968
$ws_CRP.Range("J58").Value = 3405.8
$ws_CRP.Range("K58").Value = 1606.0968
$ws_CRP.Range("L58").Value = 3405.8
$ws_CRP.Range("M58").Value = -1403.0968
$ws_CRP.Range("N58").Value = -3811.8

# CRP row 134
$ws_CRP.Range("H134").Value = 6037.36
$ws_CRP.Range("I134").Value = 6944.2354
$ws_CRP.Range("J134").Value = 4110.25
$ws_CRP.Range("K134").Value = 20832.7062
$ws_CRP.Range("L134").Value = 12330.75
$ws_CRP.Range("M134").Value = -18297.7062
$ws_CRP.Range("N134").Value = -17400.75

# CRP row 136
$ws_CRP.Range("H136").Value = 1856.0555
$ws_CRP.Range("I136").Value = 1606.0968
$ws_CRP.Range("J136").Value = 3405.8
$ws_CRP.Range("K136").Value = 4818.2904
$ws_CRP.Range("L136").Value = 10217.4
$ws_CRP.Range("M136").Value = -2268.2904
$ws_CRP.Range("N136").Value = -15317.4

# GSM row 113
$ws_GSM.Range("H113").Value = 1269.2
$ws_GSM.Range("I113").Value = 1363.0834
$ws_GSM.Range("J113").Value = 1182.5385
$ws_GSM.Range("K113").Value = 1363.0834
$ws_GSM.Range("L113").Value = 1182.5385
$ws_GSM.Range("M113").Value = 806.9166
$ws_GSM.Range("N113").Value = -5522.538500000001

# GSM row 122
$ws_GSM.Range("H122").Value = 2717.0312
$ws_GSM.Range("I122").Value = 2347.2222
$ws_GSM.Range("J122").Value = 3192.5
$ws_GSM.Range("K122").Value = 7041.6666
$ws_GSM.Range("L122").Value = 9577.5
$ws_GSM.Range("M122").Value = -4591.6666
$ws_GSM.Range("N122").Value = -14477.5

# GSM row 126
$ws_GSM.Range("H126").Value = 4155.5557
$ws_GSM.Range("I126").Value = 2957.1428
$ws_GSM.Range("J126").Value = 5446.154
$ws_GSM.Range("K126").Value = 8871.428400000001
$ws_GSM.Range("L126").Value = 16338.462
$ws_GSM.Range("M126").Value = -6401.428400000001
$ws_GSM.Range("N126").Value = -21278.462

# GSM row 127
$ws_GSM.Range("H127").Value = 20142
$ws_GSM.Range("J127").Value = 20142
$ws_GSM.Range("L127").Value = 20142
$ws_GSM.Range("N127").Value = -30062

# GSM row 132
$ws_GSM.Range("H132").Value = 3768.9443
$ws_GSM.Range("I132").Value = 2000
$ws_GSM.Range("J132").Value = 4274.357
$ws_GSM.Range("K132").Value = 6000
$ws_GSM.Range("L132").Value = 12823.071
$ws_GSM.Range("M132").Value = -3470
$ws_GSM.Range("N132").Value = -17883.071

# LTW row 61
$ws_LTW.Range("H61").Value = 1891.5264
$ws_LTW.Range("I61").Value = 1860.75
$ws_LTW.Range("J61").Value = 1944.2858
$ws_LTW.Range("K61").Value = 1860.75
$ws_LTW.Range("L61").Value = 1944.2858
$ws_LTW.Range("M61").Value = -1658.75
$ws_LTW.Range("N61").Value = -2348.2858

# LTW row 93
$ws_LTW.Range("H93").Value = 2386.6086
$ws_LTW.Range("I93").Value = 1817.1765
$ws_LTW.Range("K93").Value = 1817.1765
$ws_LTW.Range("M93").Value = -569.1765

# LTW row 113
$ws_LTW.Range("H113").Value = 1891.5264
$ws_LTW.Range("I113").Value = 1860.75
$ws_LTW.Range("J113").Value = 1944.2858
$ws_LTW.Range("K113").Value = 1860.75
$ws_LTW.Range("L113").Value = 1944.2858
$ws_LTW.Range("M113").Value = 309.25
$ws_LTW.Range("N113").Value = -6284.2858

# LTW row 132
$ws_LTW.Range("H132").Value = 4014.7192
$ws_LTW.Range("I132").Value = 1556.3478
$ws_LTW.Range("J132").Value = 5677.7354
$ws_LTW.Range("K132").Value = 4669.0434
$ws_LTW.Range("L132").Value = 17033.2062
$ws_LTW.Range("M132").Value = -2139.0434
$ws_LTW.Range("N132").Value = -22093.2062

# LTW row 136
$ws_LTW.Range("H136").Value = 4803.5
$ws_LTW.Range("I136").Value = 2386.0667
$ws_LTW.Range("J136").Value = 8100
$ws_LTW.Range("K136").Value = 7158.2001
$ws_LTW.Range("L136").Value = 24300
$ws_LTW.Range("M136").Value = -4608.2001
$ws_LTW.Range("N136").Value = -29400

# WVR row 107
$ws_WVR.Range("H107").Value = 713.5238000000001
$ws_WVR.Range("I107").Value = 604.7646999999999
$ws_WVR.Range("K107").Value = 1814.2941
$ws_WVR.Range("M107").Value = 105.7059000000002

# WVR row 122
$ws_WVR.Range("H122").Value = 2838.85
$ws_WVR.Range("I122").Value = 1860.2593
$ws_WVR.Range("K122").Value = 5580.7779
$ws_WVR.Range("M122").Value = -3130.7779

# WVR row 125
$ws_WVR.Range("H125").Value = 42043.5
$ws_WVR.Range("J125").Value = 42043.5
$ws_WVR.Range("L125").Value = 42043.5
$ws_WVR.Range("N125").Value = -51883.5

# WVR row 126
$ws_WVR.Range("H126").Value = 335454.75
$ws_WVR.Range("J126").Value = 822694.1
$ws_WVR.Range("L126").Value = 2468082.3
$ws_WVR.Range("N126").Value = -2473022.3

# WVR row 132
$ws_WVR.Range("H132").Value = 5651090.5
$ws_WVR.Range("I132").Value = 510.73685
$ws_WVR.Range("K132").Value = 1532.21055
$ws_WVR.Range("M132").Value = 997.78945
